# Update the three-digit / one-digit division problems in the table.
# Each "old=" text is unique within the document, so a simple
# Find/Replace (wdReplaceAll) per pair is safe and precise.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "195÷9="; New = "431÷7=" },
    @{ Old = "453÷8="; New = "203÷3=" },
    @{ Old = "893÷4="; New = "779÷3=" },
    @{ Old = "913÷2="; New = "605÷3=" },
    @{ Old = "952÷2="; New = "701÷3=" },
    @{ Old = "879÷5="; New = "115÷7=" },
    @{ Old = "355÷3="; New = "494÷9=" },
    @{ Old = "644÷9="; New = "733÷7=" },
    @{ Old = "406÷3="; New = "976÷5=" },
    @{ Old = "368÷5="; New = "677÷5=" },
    @{ Old = "172÷8="; New = "147÷4=" },
    @{ Old = "945÷4="; New = "350÷5=" },
    @{ Old = "401÷8="; New = "446÷5=" },
    @{ Old = "119÷4="; New = "326÷5=" },
    @{ Old = "556÷8="; New = "304÷6=" },
    @{ Old = "294÷9="; New = "593÷3=" },
    @{ Old = "899÷8="; New = "996÷9=" },
    @{ Old = "707÷4="; New = "591÷9=" },
    @{ Old = "143÷6="; New = "352÷9=" },
    @{ Old = "931÷5="; New = "636÷5=" },
    @{ Old = "726÷2="; New = "547÷6=" },
    @{ Old = "203÷2="; New = "827÷9=" },
    @{ Old = "792÷2="; New = "845÷9=" },
    @{ Old = "335÷9="; New = "573÷9=" },
    @{ Old = "667÷3="; New = "748÷6=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
